# This script updates rows in Sheet1 where the "email" (B), and in some
# cases "status" (C) / "priority" (D), columns were left blank.
# Blank emails are now explicitly marked as "N/A". Rows where email,
# status AND priority were all blank additionally get status="N/A"
# and priority=0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where only column B (email) was blank; C (status) and D (priority)
# already had values.
$bOnlyRanges = @("B14:B23", "B191:B206", "B595:B600")
foreach ($rng in $bOnlyRanges) {
    $ws.Range($rng).Value = "N/A"
}

# Rows where columns B (email), C (status) and D (priority) were all blank.
$bcdRowRanges = @(@(93, 103), @(588, 591))
foreach ($pair in $bcdRowRanges) {
    $startRow = $pair[0]
    $endRow = $pair[1]
    $ws.Range("B$startRow`:B$endRow").Value = "N/A"
    $ws.Range("C$startRow`:C$endRow").Value = "N/A"
    $ws.Range("D$startRow`:D$endRow").Value = 0
}
